$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# New row 35 values
$ws.Range("A35").Value = 1.8
$ws.Range("A35").NumberFormat = "0.00"
$ws.Range("B35").Value = 0.3
$ws.Range("C35").Value = 250
$ws.Range("D35").Value = 25
$ws.Range("E35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = "'true"
$ws.Range("I35").ClearFormats()
$ws.Range("J35").Value = 1.6
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 580
$ws.Range("M35").Value = "nach 58 M Comp => 1.0 rating mit 18-5-5-9 netz"
$ws.Range("N35").Value = "0.05 (phased)"

# Update selection to reflect the new active cell
$ws.Range("N35").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
